# Auto-generated edit script: updates crypto price/volume table
# to match the "Updated cryptos list" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look numeric but must remain plain text (matching the
# original inline-string cell type), so force Text number format first.
$textCells = @("D4","D5","D6","D7","D8","D9","D10","D11","D12","D15","D17","D20","D21","D22","D23","D24","D25","D26","D27","D28","D29","D30","D31","D32","D34","D35","D36","D37","D38","D39","D40","D41","D42","D43","D46","D47","D48","D49","D50","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the new cell values (price + 1h volume change columns).
$ws.Range("D2").Value = "26.044.83"
$ws.Range("E2").Value = "  -0.60%  "
$ws.Range("D3").Value = "1.650.47"
$ws.Range("E3").Value = "  -0.55%  "
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.19%  "
$ws.Range("D5").Value = "217.26"
$ws.Range("E5").Value = "  +0.03%  "
$ws.Range("D6").Value = "0.5260"
$ws.Range("E6").Value = "  +1.25%  "
$ws.Range("D7").Value = "1.002"
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("D8").Value = "0.2593"
$ws.Range("E8").Value = "  -1.70%  "
$ws.Range("D9").Value = "0.06306"
$ws.Range("E9").Value = "  +0.57%  "
$ws.Range("D10").Value = "20.30"
$ws.Range("E10").Value = "  -2.21%  "
$ws.Range("D11").Value = "0.07793"
$ws.Range("E11").Value = "  -0.02%  "
$ws.Range("D12").Value = "4.511"
$ws.Range("E12").Value = "  +0.94%  "
$ws.Range("D13").Value = "1.569.63"
$ws.Range("E13").Value = "  -5.40%  "
$ws.Range("D14").Value = "1.877.50"
$ws.Range("E14").Value = "  -0.51%  "
$ws.Range("D15").Value = "0.5479"
$ws.Range("E15").Value = "  +0.16%  "
$ws.Range("D16").Value = "0.0₅8185"
$ws.Range("E16").Value = "  +0.80%  "
$ws.Range("D17").Value = "65.46"
$ws.Range("E17").Value = "  +0.72%  "
$ws.Range("D18").Value = "26.054.83"
$ws.Range("E18").Value = "  -0.56%  "
$ws.Range("E19").Value = "  -0.13%  "
$ws.Range("D20").Value = "4.567"
$ws.Range("E20").Value = "  -0.98%  "
$ws.Range("D21").Value = "190.37"
$ws.Range("E21").Value = "  -0.83%  "
$ws.Range("D22").Value = "10.08"
$ws.Range("E22").Value = "  +0.23%  "
$ws.Range("D23").Value = "6.006"
$ws.Range("E23").Value = "  +0.09%  "
$ws.Range("D24").Value = "1.003"
$ws.Range("E24").Value = "  -0.20%  "
$ws.Range("D25").Value = "143.16"
$ws.Range("E25").Value = "  +2.90%  "
$ws.Range("D26").Value = "0.1233"
$ws.Range("E26").Value = "  +0.90%  "
$ws.Range("D27").Value = "7.215"
$ws.Range("E27").Value = "  -0.98%  "
$ws.Range("D28").Value = "16.00"
$ws.Range("E28").Value = "  -0.92%  "
$ws.Range("D29").Value = "1.434"
$ws.Range("E29").Value = "  -0.04%  "
$ws.Range("D30").Value = "0.05814"
$ws.Range("E30").Value = "  -2.24%  "
$ws.Range("D31").Value = "1.269"
$ws.Range("E31").Value = "  -0.56%  "
$ws.Range("D32").Value = "3.542"
$ws.Range("E32").Value = "  -0.23%  "
$ws.Range("E33").Value = "  -0.33%  "
$ws.Range("D34").Value = "1.588"
$ws.Range("E34").Value = "  +0.69%  "
$ws.Range("D35").Value = "2.795"
$ws.Range("E35").Value = "  +1.00%  "
$ws.Range("D36").Value = "2.415"
$ws.Range("E36").Value = "  -0.21%  "
$ws.Range("D37").Value = "0.9413"
$ws.Range("E37").Value = "  -2.01%  "
$ws.Range("D38").Value = "0.5750"
$ws.Range("E38").Value = "  +1.12%  "
$ws.Range("D39").Value = "0.01604"
$ws.Range("D40").Value = "0.8472"
$ws.Range("E40").Value = "  -0.44%  "
$ws.Range("D41").Value = "104.67"
$ws.Range("E41").Value = "  +4.33%  "
$ws.Range("D42").Value = "1.002"
$ws.Range("E42").Value = "  -0.10%  "
$ws.Range("D43").Value = "5.705"
$ws.Range("E43").Value = "  -4.94%  "
$ws.Range("D44").Value = "1.029.57"
$ws.Range("E44").Value = "  +2.62%  "
$ws.Range("D45").Value = "1.792.02"
$ws.Range("E45").Value = "  -0.53%  "
$ws.Range("D46").Value = "57.04"
$ws.Range("E46").Value = "  +1.00%  "
$ws.Range("D47").Value = "0.9997"
$ws.Range("E47").Value = "  -0.19%  "
$ws.Range("D48").Value = "0.4327"
$ws.Range("E48").Value = "  +0.04%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "7.835"
$ws.Range("E49").Value = "  -2.19%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "0.05142"
$ws.Range("E50").Value = "  -0.35%  "
$ws.Range("D51").Value = "1.451"
$ws.Range("E51").Value = "  -0.13%  "

Write-Output "Updated cryptos list"
